# [NEAT-671] Add GUID to source system
# Adds a new "classicSourceSystemGUID" property to the ClassicSourceSystem
# mapping (Properties sheet) and registers the corresponding container
# entry (Containers sheet).

$wb = $excel.ActiveWorkbook

# --- Properties sheet -----------------------------------------------------
# Insert a new row right after the existing ClassicSourceSystem "name"
# property row, describing the new classicSourceSystemGUID property.
$wsProps = $wb.Worksheets.Item("Properties")
$wsProps.Rows.Item(54).Insert()

$wsProps.Cells.Item(54, 1).Value2 = "ClassicSourceSystem"       # Class
$wsProps.Cells.Item(54, 2).Value2 = "classicSourceSystemGUID"   # Property
$wsProps.Cells.Item(54, 6).Value2 = "text"                      # Value Type
$wsProps.Cells.Item(54, 7).Value2 = $true                       # Nullable
$wsProps.Cells.Item(54, 8).Value2 = $false                      # Immutable
$wsProps.Cells.Item(54, 9).Value2 = $false                      # Is List
$wsProps.Cells.Item(54, 11).Value2 = "ClassicSourceSystem"      # Container
$wsProps.Cells.Item(54, 12).Value2 = "classicSourceSystemGUID"  # Container Property

# --- Containers sheet -------------------------------------------------------
# Register the classicSourceSystemGUID container/property as a plain node.
$wsContainers = $wb.Worksheets.Item("Containers")
$wsContainers.Cells.Item(4, 1).Value2 = "ClassicSourceSystem"   # Container
$wsContainers.Cells.Item(4, 5).Value2 = "node"                  # Used For

# --- Restore the cursor/selection positions left behind by the edit --------
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Activate()
$wsViews.Range("C13").Select()

$wsProps.Activate()
$wsProps.Range("E78").Select()
